# Mexico Liga TDP - "Atualização de bases das ligas, do dia: 24-02-2024 às 23:13"
#
# The update swaps the full data payload (columns B:AC) between pairs of rows
# (and one 3-way rotation), which is how this source workbook de-duplicates /
# re-orders same-day fixtures after a re-scrape. Column A (the running row
# index) is left untouched since it always matches the row's position.
# One row (696) also receives a standalone odds correction (no swap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($sheet, [int]$rowA, [int]$rowB)
    $ra = $sheet.Range("B$rowA`:AC$rowA")
    $rb = $sheet.Range("B$rowB`:AC$rowB")
    $tmp = $ra.Value2
    $ra.Value2 = $rb.Value2
    $rb.Value2 = $tmp
}

function Rotate-Rows {
    # new(row1) = old(row2); new(row2) = old(row3); new(row3) = old(row1)
    param($sheet, [int]$row1, [int]$row2, [int]$row3)
    $r1 = $sheet.Range("B$row1`:AC$row1")
    $r2 = $sheet.Range("B$row2`:AC$row2")
    $r3 = $sheet.Range("B$row3`:AC$row3")
    $scratch = $sheet.Range("B1000:AC1000")
    $scratch.Value2 = $r1.Value2
    $r1.Value2 = $r2.Value2
    $r2.Value2 = $r3.Value2
    $r3.Value2 = $scratch.Value2
    $scratch.ClearContents()
}

# Pairs of rows whose full record (B:AC) swapped places.
Swap-Rows $ws 579 581
Swap-Rows $ws 598 599
Swap-Rows $ws 633 634
Swap-Rows $ws 647 648
Swap-Rows $ws 664 665

# Three rows whose records cyclically rotated.
Rotate-Rows $ws 636 637 638

# Standalone closing-odds correction on row 696 (N,O,P,Q,R,S,U,V only; T unchanged).
$ws.Range("N696").Value2 = 6.5
$ws.Range("O696").Value2 = 4.75
$ws.Range("P696").Value2 = 1.363
$ws.Range("Q696").Value2 = 1.5
$ws.Range("R696").Value2 = 1.9
$ws.Range("S696").Value2 = 1.9
$ws.Range("U696").Value2 = 1.925
$ws.Range("V696").Value2 = 1.875
